$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (shifts the existing data rows down by one).
$ws.Rows.Item(2).Insert()

# The insert copies the bold/centered header formatting down onto the new
# row - clear that back to the plain default formatting used by the other
# data rows (this is the alignment glitch being fixed).
$ws.Rows.Item(2).ClearFormats()

# Populate the new row with its data.
$ws.Cells.Item(2, 1).Value = "Move to location (9, 5) and remove the toolkit."
$ws.Cells.Item(2, 2).Value = 75.906284
$ws.Cells.Item(2, 3).Value = 10794
# Cost column is stored as text in this sheet, so force text (leading
# apostrophe) instead of letting it be auto-detected as a number.
$ws.Cells.Item(2, 4).Value = "'0.03354"

# The sheet keeps a fixed 10 data rows, so the row that gets pushed past
# the bottom (old row 11) is dropped.
$ws.Rows.Item(12).Delete()
